# Applies the weekly Fruta/Hortaliza data refresh for the
# "Macroferia Regional de Talca - Arveja Verde" sheet:
#  - rows 30-43 get updated field values (prices/dates/origins/etc.)
#  - two new rows (44, 45) are appended with new records
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 30-43 with changed cell values ---
# Row 30
$ws.Range("D30").Value = "2021-10-19"
$ws.Range("J30").Value = 200
$ws.Range("O30").Value = 'Región de O''Higgins'

# Row 31
$ws.Range("D31").Value = "2021-10-19"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 19000
$ws.Range("L31").Value = 19000
$ws.Range("M31").Value = 19000
$ws.Range("P31").Value = 760

# Row 32
$ws.Range("D32").Value = "2020-12-03"
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 20000
$ws.Range("O32").Value = 'Región de La Araucanía'
$ws.Range("P32").Value = 800

# Row 33
$ws.Range("D33").Value = "2020-11-26"
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 16000
$ws.Range("L33").Value = 16000
$ws.Range("M33").Value = 16000
$ws.Range("O33").Value = 'Región del Maule'
$ws.Range("P33").Value = 640

# Row 34
$ws.Range("D34").Value = "2020-11-30"
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 15000
$ws.Range("N34").Value = '$/saco 25 kilos'
$ws.Range("O34").Value = 'Región del Maule'
$ws.Range("P34").Value = 600
$ws.Range("Q34").Value = 25

# Row 35
$ws.Range("D35").Value = "2021-03-03"
$ws.Range("N35").Value = '$/saco 25 kilos'
$ws.Range("P35").Value = 880
$ws.Range("Q35").Value = 25

# Row 36
$ws.Range("D36").Value = "2021-02-25"
$ws.Range("H36").Value = 'Sin especificar'
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 22000
$ws.Range("L36").Value = 22000
$ws.Range("M36").Value = 22000
$ws.Range("N36").Value = '$/saco 30 kilos'
$ws.Range("O36").Value = 'Región de La Araucanía'
$ws.Range("P36").Value = 22000
$ws.Range("Q36").Value = 1

# Row 37
$ws.Range("D37").Value = "2021-02-16"

# Row 38
$ws.Range("D38").Value = "2021-09-14"
$ws.Range("H38").Value = 'Perfection'
$ws.Range("J38").Value = 150
$ws.Range("K38").Value = 35000
$ws.Range("L38").Value = 35000
$ws.Range("M38").Value = 35000
$ws.Range("N38").Value = '$/malla 25 kilos'
$ws.Range("O38").Value = 'Provincia del Elquí'
$ws.Range("P38").Value = 1400

# Row 39
$ws.Range("D39").Value = "2021-02-15"
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 22000
$ws.Range("L39").Value = 22000
$ws.Range("M39").Value = 22000
$ws.Range("N39").Value = '$/saco 30 kilos'
$ws.Range("O39").Value = 'Región de La Araucanía'
$ws.Range("P39").Value = 22000
$ws.Range("Q39").Value = 1

# Row 40
$ws.Range("D40").Value = "2020-11-24"
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 17000
$ws.Range("L40").Value = 18000
$ws.Range("M40").Value = 17500
$ws.Range("O40").Value = 'Región del Maule'
$ws.Range("P40").Value = 700

# Row 41
$ws.Range("D41").Value = "2020-12-01"
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 16000
$ws.Range("L41").Value = 16000
$ws.Range("M41").Value = 16000
$ws.Range("O41").Value = 'Región del Maule'
$ws.Range("P41").Value = 640

# Row 42
$ws.Range("D42").Value = "2021-03-08"
$ws.Range("J42").Value = 300
$ws.Range("K42").Value = 22000
$ws.Range("L42").Value = 22000
$ws.Range("M42").Value = 22000
$ws.Range("O42").Value = 'Región de La Araucanía'
$ws.Range("P42").Value = 880

# Row 43
$ws.Range("D43").Value = "2020-12-07"
$ws.Range("J43").Value = 250
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = 20000
$ws.Range("P43").Value = 800

# --- Add new rows 44-45 (appended records) ---
# Row 44
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A44").Value = 5
$ws.Range("B44").Value = 'Macroferia Regional de Talca'
$ws.Range("C44").Value = 'Maule'
$ws.Range("D44").Value = "2021-03-15"
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = 100112022
$ws.Range("G44").Value = 'Arveja Verde'
$ws.Range("H44").Value = 'Perfection'
$ws.Range("I44").Value = 'Primera'
$ws.Range("J44").Value = 250
$ws.Range("K44").Value = 23000
$ws.Range("L44").Value = 23000
$ws.Range("M44").Value = 23000
$ws.Range("N44").Value = '$/saco 25 kilos'
$ws.Range("O44").Value = 'Región del Maule'
$ws.Range("P44").Value = 920
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = 'Hortaliza'

# Row 45
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = 'Macroferia Regional de Talca'
$ws.Range("C45").Value = 'Maule'
$ws.Range("D45").Value = "2021-02-23"
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = 100112022
$ws.Range("G45").Value = 'Arveja Verde'
$ws.Range("H45").Value = 'Sin especificar'
$ws.Range("I45").Value = 'Primera'
$ws.Range("J45").Value = 200
$ws.Range("K45").Value = 22000
$ws.Range("L45").Value = 22000
$ws.Range("M45").Value = 22000
$ws.Range("N45").Value = '$/saco 25 kilos'
$ws.Range("O45").Value = 'Región de La Araucanía'
$ws.Range("P45").Value = 880
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = 'Hortaliza'
